$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: B36 was stored as text "3"; fix it to be a real number 3
$ws.Range("B36").Value = 3

# Append new row 37 with the annotation data that used to live in row 36's
# "B" slot (kept as text, matching the original quirky text-typed "3"),
# plus the new annotation columns.
$ws.Range("A37").Value = "Ruilin"
$ws.Range("B37").Value = "'3"
$ws.Range("C37").Value = "无"
$ws.Range("D37").Value = "DFT"
$ws.Range("E37").Value = "MET"
$ws.Range("F37").Value = "0c8a854c-e7df-48dd-93a0-b6771319a745"
$ws.Range("G37").Value = "H1Ww66x0-_annotated.xlsx"
$ws.Range("H37").Value = "- the proposed approach to maintain the budget is simplistic"
